$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New symbol-list values scraped on Sat Jan 14 15:48:48 UTC 2023.
# Cells are plain text (Price / Volume(1h) are formatted strings, not
# numbers), so we force Text number-format before writing, then restore
# the cell's style so no stray formatting is introduced.

$targetCells = @("D2","E2","D3","E3","D4","E4","D5","E5","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","D26","E26","D27","E27","D28","E28","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "304.28"
$ws.Range("E2").Value = "6.02%"
$ws.Range("D3").Value = "32.02"
$ws.Range("E3").Value = "8.66%"
$ws.Range("D4").Value = "5.308"
$ws.Range("E4").Value = "4.17%"
$ws.Range("D5").Value = "0.07475"
$ws.Range("E5").Value = "6.89%"
$ws.Range("E6").Value = "5.38%"
$ws.Range("D7").Value = "3.825"
$ws.Range("E7").Value = "7.99%"
$ws.Range("D8").Value = "1.470"
$ws.Range("E8").Value = "6.28%"
$ws.Range("D9").Value = "0.9217"
$ws.Range("E9").Value = "2.14%"
$ws.Range("D10").Value = "0.01755"
$ws.Range("E10").Value = "2,610.24%"
$ws.Range("D11").Value = "0.1698"
$ws.Range("E11").Value = "6.69%"
$ws.Range("D12").Value = "0.07671"
$ws.Range("E12").Value = "10.09%"
$ws.Range("D13").Value = "0.08030"
$ws.Range("E13").Value = "4.68%"
$ws.Range("D14").Value = "0.03057"
$ws.Range("E14").Value = "4.87%"
$ws.Range("D15").Value = "0.09915"
$ws.Range("E15").Value = "10.28%"
$ws.Range("D16").Value = "0.001494"
$ws.Range("E16").Value = "-5.87%"
$ws.Range("D17").Value = "0.04571"
$ws.Range("E17").Value = "1.43%"
$ws.Range("D18").Value = "0.006468"
$ws.Range("E18").Value = "1.74%"
$ws.Range("D19").Value = "3.483"
$ws.Range("E19").Value = "0.38%"
$ws.Range("E20").Value = "0.02%"
$ws.Range("D21").Value = "0.3335"
$ws.Range("E21").Value = "3.15%"
$ws.Range("D22").Value = "0.1346"
$ws.Range("E22").Value = "1.64%"
$ws.Range("D23").Value = "4.510"
$ws.Range("E23").Value = "12.19%"
$ws.Range("D24").Value = "0.1624"
$ws.Range("E24").Value = "4.50%"
$ws.Range("D25").Value = "0.001219"
$ws.Range("D26").Value = "0.004416"
$ws.Range("E26").Value = "0.84%"
$ws.Range("D27").Value = "0.0001403"
$ws.Range("E27").Value = "20.37%"
$ws.Range("D28").Value = "0.0001744"
$ws.Range("E28").Value = "8.18%"
$ws.Range("D40").Value = "0.04528"
$ws.Range("E40").Value = "5.57%"
$ws.Range("D41").Value = "0.007214"
$ws.Range("E41").Value = "5.70%"
$ws.Range("D42").Value = "0.1342"
$ws.Range("E42").Value = "7.77%"
$ws.Range("D43").Value = "0.002205"
$ws.Range("E43").Value = "1.05%"
$ws.Range("D44").Value = "0.01268"
$ws.Range("E44").Value = "9.95%"
$ws.Range("D45").Value = "0.00006151"
$ws.Range("E45").Value = "7.00%"
$ws.Range("D46").Value = "0.7092"
$ws.Range("E46").Value = "-63.24%"
$ws.Range("D47").Value = "0.01302"
$ws.Range("E47").Value = "0.00%"

foreach ($addr in $targetCells) {
    $ws.Range($addr).Style = "Normal"
}
